# Insert a new daily price record at row 93 for "Berenjena" (Vega Modelo de
# Temuco), pushing the existing rows 93:189 down to 94:190.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 93:189 down by one row (creates a new blank row 93).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new record.
$ws.Range("A93").Value = 10
$ws.Range("B93").Value = "Vega Modelo de Temuco"
$ws.Range("C93").Value = "La Araucanía"
$ws.Range("D93").Value = 44484
$ws.Range("E93").Value = 9
$ws.Range("F93").Value = 100112001
$ws.Range("G93").Value = "Berenjena"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 50
$ws.Range("K93").Value = 10000
$ws.Range("L93").Value = 10000
$ws.Range("M93").Value = 10000
$ws.Range("N93").Value = "`$/caja 60 unidades"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 167
$ws.Range("Q93").Value = 60
$ws.Range("R93").Value = "Hortaliza"
